# Rework the DivineDomains sheet: collapse the wide "one column per domain"
# header row down to two domains (Ambition, Arcana) and add the supporting
# data rows (subclass name, bonus-spell links, and feature lists) underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused domain headers in C1:W1 (Blood Domain .. Zeal Domain).
# Clear() (not ClearContents()) drops the cell nodes entirely instead of
# leaving behind empty-but-present <c> elements.
$ws.Range("C1:W1").Clear()

# Row 1 stays as the two remaining domain headers.
$ws.Cells.Item(1, 1).Value = "Ambition Domain"
$ws.Cells.Item(1, 2).Value = "Arcana Domain"

# Row 2-4: "None" placeholder in both columns.
$ws.Cells.Item(2, 1).Value = "None"
$ws.Cells.Item(2, 2).Value = "None"

$ws.Cells.Item(3, 1).Value = "None"
$ws.Cells.Item(3, 2).Value = "None"

$ws.Cells.Item(4, 1).Value = "None"
$ws.Cells.Item(4, 2).Value = "None"

# Row 5: "None" / "Arcana".
$ws.Cells.Item(5, 1).Value = "None"
$ws.Cells.Item(5, 2).Value = "Arcana"

# Row 6: bonus-spell workbook links.
$ws.Cells.Item(6, 1).Value = "true=classes/cleric/domains/BonusSpellsAmbitionDomain.xlsx"
$ws.Cells.Item(6, 2).Value = "true=classes/cleric/domains/BonusSpellsArcaneDomain.xlsx"

# Row 7: "false=" marker in both columns.
$ws.Cells.Item(7, 1).Value = "false="
$ws.Cells.Item(7, 2).Value = "false="

# Row 8: level/feature lists for each domain.
$ws.Cells.Item(8, 1).Value = "1/Warding Flare=2/Channel Divinity: Invoke Duplicity=6/Channel Divinity: Cloak of Shadows=8/Potent Spellcasting=17/Improved Duplicity"
$ws.Cells.Item(8, 2).Value = "2/Channel Divinity: Arcane Abjuration=6/Spell Breaker=8/Potent Spellcasting"

# Match row 1's cell style (s="1") on the new data rows, same as the rest
# of the header row.
$ws.Range("A1").Copy()
$ws.Range("A2:B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
